# Update "handback-status" timestamps generated by the handback report.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for c7f153eb-...md
$wsOverview.Range("G3").Value = "2016-08-18 16:44:35"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for c7f153eb-...
$wsZhCn.Range("H3").Value = "2016-08-18 16:44:30"
$wsZhCn.Range("K3").Value = "2016-08-18 16:44:48"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for c7f153eb-...
$wsDeDe.Range("H3").Value = "2016-08-18 16:44:35"
$wsDeDe.Range("K3").Value = "2016-08-18 16:44:55"
